$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (Generation 0)
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4878.515958786011
$ws.Range("F2").Value = 92
$ws.Range("G2").Value = 16

# Add new row 3 (Generation 1)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "wins_data"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 4468.364238739014
$ws.Range("F3").Value = 91
$ws.Range("G3").Value = 16

# Add new row 4 (Generation 2)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "wins_data"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 4360.931396484375
$ws.Range("F4").Value = 87
$ws.Range("G4").Value = 15

# Add new row 5 (Generation 3)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "wins_data"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 4282.888650894165
$ws.Range("F5").Value = 89
$ws.Range("G5").Value = 16
